# Update countries & provincias Spain
# Refresh of the COVID-19 "Pais" data table: some country stats changed,
# which moves a few tied rows around (the data stays sorted by total
# cases, column B, descending), and the "last updated" footer timestamp
# advances from 06:40 to 07:10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Uzbekistan (row 78): totals/new-cases/active-cases refreshed ---
$ws.Range("B78").Value = 3488
$ws.Range("C78").Value = 20
$ws.Range("E78").Value = 746

# --- Tailandia (row 80): totals/new-cases/recovered/active refreshed ---
$ws.Range("B80").Value = 3077
$ws.Range("C80").Value = 1
$ws.Range("D80").Value = 2961
$ws.Range("E80").Value = 59

# --- El Salvador (row 88): recovered/active/critical/deaths refreshed ---
$ws.Range("D88").Value = 1017
$ws.Range("E88").Value = 1218
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 43

# --- Curazao / Fiyi (rows 198-199) swap places (tied on column B) ---
$ws.Range("A198").Value = "Fiyi"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# --- Montserrat / Seychelles (rows 210-211) swap places (tied on column B) ---
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- Islas Virgenes Britanicas / Papua Nueva Guinea (rows 213-214) swap places ---
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 07:10"
